# Apply corrected Diebold-Mariano values to P_valores and Estadisticos_DM sheets

$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores" ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.7411977073101284
$wsP.Range("D2").Value = 0.6953172266069783
$wsP.Range("E2").Value = 0.691196157911945
$wsP.Range("F2").Value = 0.8727831780625785

$wsP.Range("B3").Value = 0.7411977073101284
$wsP.Range("D3").Value = 0.9700568259076627
$wsP.Range("E3").Value = 0.8371838364103432
$wsP.Range("F3").Value = 0.6648893101372944

$wsP.Range("B4").Value = 0.6953172266069783
$wsP.Range("C4").Value = 0.9700568259076627
$wsP.Range("E4").Value = 0.8015557060228455
$wsP.Range("F4").Value = 0.6912865039246854

$wsP.Range("B5").Value = 0.691196157911945
$wsP.Range("C5").Value = 0.8371838364103432
$wsP.Range("D5").Value = 0.8015557060228455
$wsP.Range("F5").Value = 0.848709224014087

$wsP.Range("B6").Value = 0.8727831780625785
$wsP.Range("C6").Value = 0.6648893101372944
$wsP.Range("D6").Value = 0.6912865039246854
$wsP.Range("E6").Value = 0.848709224014087

# --- Sheet "Estadisticos_DM" ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 0.3344631918642363
$wsE.Range("D2").Value = 0.3968285104233975
$wsE.Range("E2").Value = 0.4025057226320155
$wsE.Range("F2").Value = 0.1620017124876339

$wsE.Range("B3").Value = -0.3344631918642363
$wsE.Range("D3").Value = -0.03796647776099236
$wsE.Range("E3").Value = -0.2079463849863685
$wsE.Range("F3").Value = -0.4390751520159932

$wsE.Range("B4").Value = -0.3968285104233975
$wsE.Range("C4").Value = 0.03796647776099236
$wsE.Range("E4").Value = -0.2543916657058904
$wsE.Range("F4").Value = -0.4023811175816107

$wsE.Range("B5").Value = -0.4025057226320155
$wsE.Range("C5").Value = 0.2079463849863685
$wsE.Range("D5").Value = 0.2543916657058904
$wsE.Range("F5").Value = -0.1930263597917401

$wsE.Range("B6").Value = -0.1620017124876339
$wsE.Range("C6").Value = 0.4390751520159932
$wsE.Range("D6").Value = 0.4023811175816107
$wsE.Range("E6").Value = 0.1930263597917401
